# Apply "contingencies with rene fine" edit:
# - Extend the table by two columns (P, Q) with header values 14, 15 on row 1
# - For rows 2..25: swap I<->K style of change (I:1->2, K:2->1, M:1->2, O:2->1)
#   and add new P/Q columns with value 2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new header cells P1, Q1 ---
# Copy formatting from O1 (bold, centered, bordered) onto the new header cells
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Rows 2-25: update existing values and add new P/Q columns ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}
